# Select the "Granty_złożone" sheet and remove the "Razem" (total) summary
# rows for 2021, 2020 and 2019 from it, then leave the active cell on A19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Granty_złożone")
$ws.Activate()

# Delete entire rows from the bottom up so earlier row numbers stay valid.
$ws.Rows.Item(55).Delete()
$ws.Rows.Item(37).Delete()
$ws.Rows.Item(19).Delete()

$ws.Rows.Item(19).Select()
